# Actualización automática del tracker
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the result of the match in row 73 (now resolved as a loss/"Fallo")
$ws.Range("G73").Value = "Fallo"
$ws.Range("H73").Value = -1

# Append newly tracked matches as rows 85 and 86.
# "fecha" is stored as plain text (e.g. "2025-08-06"), not a date serial,
# so the value is entered via a formula and then flattened to a static
# value (copy / paste-values) to avoid Excel's automatic date detection
# while keeping the cell's style untouched (no quote-prefix formatting).
$ws.Range("A85").Value = 14350906
$ws.Range("B85").Formula = "=TEXT(DATE(2025,8,6),""yyyy-mm-dd"")"
$ws.Range("B85").Copy()
$ws.Range("B85").PasteSpecial(-4163)
$ws.Range("C85").Value = "Mariano Navone"
$ws.Range("D85").Value = "Jan-Lennard Struff"
$ws.Range("E85").Value = "Gana Mariano Navone"
$ws.Range("F85").Value = 2.5

$ws.Range("A86").Value = 14358494
$ws.Range("B86").Formula = "=TEXT(DATE(2025,8,6),""yyyy-mm-dd"")"
$ws.Range("B86").Copy()
$ws.Range("B86").PasteSpecial(-4163)
$ws.Range("C86").Value = "Adrian Mannarino"
$ws.Range("D86").Value = "Dalibor Svrcina"
$ws.Range("E86").Value = "Gana Dalibor Svrcina"
$ws.Range("F86").Value = 1.83
